# DLP woordenlijst - add new Court Position / Lifestyle / Education terms
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 38-41 (Engelse term / Nederlandse term)
$ws.Range("A38").Value = "Grant Title"
$ws.Range("B38").Value = "Verleen Titel"

$ws.Range("A39").Value = "Hire (Court Position)"
$ws.Range("B39").Value = "Benoem"

$ws.Range("A40").Value = "Lifestyle"
$ws.Range("B40").Value = "Levensstijl"

$ws.Range("A41").Value = "Education Focus"
$ws.Range("B41").Value = "Onderwijzingsfocus"

# Grow the "Tabel1" table so the new rows are included
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C41"))

# Scroll/select like the saved workbook (cursor moved to the next empty row)
$win = $excel.ActiveWindow
$ws.Range("A42").Select()
$win.ScrollRow = 19
$win.ScrollColumn = 1
